# "suite de la jframe"
# Split the run containing "système) -> un paragraphe pour détailler ..."
# into three runs, inserting a new run with the text
# "(voir diagramme de Classe en MOO)" right after "système) " and before
# the "->" that follows it.
#
# Using a plain Find/Replace (or Range.InsertAfter/InsertBefore/Text=) at
# that spot re-flows and coalesces the whole run of text into a single
# <w:r>, which loses the three-way run split seen in the target document.
# Assigning a Range's .FormattedText to another (collapsed) Range instead
# performs a clean structural insert without touching the formatting/text
# of the surrounding runs, so we build the new text in a disposable
# scratch paragraph, grab its FormattedText, splice it in at the right
# spot, then remove the scratch paragraph again.

$d = $word.ActiveDocument

$insertText = "(voir diagramme de Classe en MOO)"
$insertLen = $insertText.Length

# --- 1. Create a throwaway paragraph at the end of the document and type
#        the new text into it, purely so we have a Range whose
#        FormattedText we can copy from. ---
$endPos = $d.Content.End
$scratchAnchor = $d.Range($endPos, $endPos)
$scratchAnchor.InsertParagraphAfter()

$scratchStart = $endPos
$scratchRange = $d.Range($scratchStart, $scratchStart)
$scratchRange.InsertBefore($insertText)

$newTextRange = $d.Range($scratchStart, $scratchStart + $insertLen)
$newFormattedText = $newTextRange.FormattedText

# --- 2. Find the real insertion point: right after "système) " and
#        right before the "-> un paragraphe ..." that follows it. ---
$full = $d.Content.Text
$anchorText = "système) "
$anchorIndex = $full.IndexOf($anchorText)
$targetPos = $anchorIndex + $anchorText.Length

$targetRange = $d.Range($targetPos, $targetPos)
$targetRange.FormattedText = $newFormattedText

# --- 3. Clean up: delete the scratch paragraph (text + paragraph mark)
#        again; its position has shifted forward by however much text we
#        just spliced in at the earlier target position. ---
$full2 = $d.Content.Text
$scratchStart2 = $scratchStart + $insertLen
$cleanupRange = $d.Range($scratchStart2, $full2.Length)
$cleanupRange.Delete()
